$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix bug related to selection of successful experiments in batch2
# Apply corrected values to the affected cells in rows 2-21
$ws.Range("E2").Value = 1.5
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 55.00000000000001
$ws.Range("D3").Value = 5
$ws.Range("F3").Value = 125
$ws.Range("G3").Value = 30
$ws.Range("H3").Value = 29
$ws.Range("C4").Value = 55.00000000000001
$ws.Range("D4").Value = 5
$ws.Range("E4").Value = 1.5
$ws.Range("G4").Value = 30
$ws.Range("H4").Value = 29
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 85.00000000000001
$ws.Range("D5").Value = 10
$ws.Range("F5").Value = 130
$ws.Range("H5").Value = 21
$ws.Range("E6").Value = 1.5
$ws.Range("H6").Value = 23
$ws.Range("D7").Value = 15
$ws.Range("F7").Value = 125
$ws.Range("H7").Value = 17
$ws.Range("E8").Value = 1.5
$ws.Range("C9").Value = 85.00000000000001
$ws.Range("E9").Value = 1.4
$ws.Range("F9").Value = 125
$ws.Range("H9").Value = 21
$ws.Range("C10").Value = 90
$ws.Range("E10").Value = 1.4
$ws.Range("H10").Value = 19
$ws.Range("C11").Value = 90
$ws.Range("F11").Value = 125
$ws.Range("H11").Value = 19
$ws.Range("D12").Value = 15
$ws.Range("E12").Value = 1.4
$ws.Range("H12").Value = 21
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 85.00000000000001
$ws.Range("D13").Value = 10
$ws.Range("G13").Value = 20
$ws.Range("H13").Value = 21
$ws.Range("F14").Value = 125
$ws.Range("H14").Value = 21
$ws.Range("C15").Value = 80
$ws.Range("F15").Value = 135
$ws.Range("H15").Value = 21
$ws.Range("C16").Value = 90
$ws.Range("D16").Value = 10
$ws.Range("F16").Value = 130
$ws.Range("H16").Value = 21
$ws.Range("C17").Value = 80
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = 1.5
$ws.Range("H17").Value = 21
$ws.Range("C18").Value = 85.00000000000001
$ws.Range("H18").Value = 21
$ws.Range("C19").Value = 85.00000000000001
$ws.Range("E19").Value = 1.5
$ws.Range("F19").Value = 130
$ws.Range("C20").Value = 85.00000000000001
$ws.Range("E20").Value = 1.6
$ws.Range("F20").Value = 130
$ws.Range("H20").Value = 21
$ws.Range("C21").Value = 80
$ws.Range("E21").Value = 1.4
$ws.Range("F21").Value = 125
